# "completed job,user for admin ui"
# Adds a new "address" column (C) to the export-jobs sheet, between
# "location" and "salary", shifting the remaining columns (salary..skills)
# one position to the right, and fills it in with per-row street addresses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; existing C:K (salary..skills) become D:L.
$null = $ws.Columns.Item(3).Insert()

# Header
$ws.Cells.Item(1, 3).Value = "address"

# Per-row address values
$ws.Cells.Item(2, 3).Value = "374 pham van dong"
$ws.Cells.Item(3, 3).Value = "375 pham van dong"
$ws.Cells.Item(4, 3).Value = "376 pham van dong"
$ws.Cells.Item(5, 3).Value = "377 pham van dong"
$ws.Cells.Item(6, 3).Value = "378 pham van dong"
$ws.Cells.Item(7, 3).Value = "379 pham van dong"
$ws.Cells.Item(8, 3).Value = "380 pham van dong"
$ws.Cells.Item(9, 3).Value = "381 pham van dong"
$ws.Cells.Item(10, 3).Value = "382 pham van dong"
$ws.Cells.Item(11, 3).Value = "383 pham van dong"

# Match the author's last selection in the workbook
$null = $ws.Range("N16").Select()
